$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at N ("Lukuvuosimaksuvelvollisuus") shifting the
# existing N..X columns (and their data) one position to the right (O..Y).
$ws.Columns("N").Insert()

# New header cell for the inserted column.
$ws.Range("N1").Value = "EI"

# Put the cursor/selection on the new column's first data cell, and scroll
# the view back to the top-left.
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("N2").Select()
